$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Monthly data refresh: revised values for the most recent existing rows
$ws.Range("B218").Value = 1827
$ws.Range("B219").Value = 1560
$ws.Range("B220").Value = 879
$ws.Range("B221").Value = 2174
$ws.Range("B222").Value = 928
$ws.Range("B223").Value = 1072

# Append the new data point for the next period: 01-07-2021 / 604
# The date-like label has to stay as plain text (like every other "Serie"
# entry in column A), but assigning the literal string directly makes Excel
# auto-convert it to a date serial number. Building it via a text formula
# forces a string result, and then copy / paste-special-values collapses the
# formula back down to a literal shared-string value without leaving any
# extra number-format/style behind.
$ws.Range("A224").Formula = "=""01-07-2021"""
$ws.Range("A224").Copy() | Out-Null
$ws.Range("A224").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B224").Value = 604
